$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '  -2.16%  '

$ws.Range("D3").Value = '2.328.10'
$ws.Range("E3").Value = '  -4.62%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").Value = '542.10'
$ws.Range("E5").Value = '  -0.90%  '

$ws.Range("D6").Value = '135.69'
$ws.Range("E6").Value = '  -6.79%  '

$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").Value = '0.518'
$ws.Range("E8").Value = '  -11.01%  '

$ws.Range("E9").Value = '  -4.79%  '

$ws.Range("E10").Value = '  -2.06%  '

$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("E12").Value = '  -2.96%  '

$ws.Range("D13").Value = '0.337'
$ws.Range("E13").Value = '  -2.90%  '

$ws.Range("D14").Value = '24.29'
$ws.Range("E14").Value = '  -5.77%  '

$ws.Range("D15").Value = '2.748.77'
$ws.Range("E15").Value = '  -4.68%  '

$ws.Range("D16").Value = '60.177.27'
$ws.Range("E16").Value = '  -1.93%  '

$ws.Range("D17").Value = '0.0000159'
$ws.Range("E17").Value = '  -5.41%  '

$ws.Range("D18").Value = '2.329.74'
$ws.Range("E18").Value = '  -4.51%  '

$ws.Range("D19").Value = '10.47'
$ws.Range("E19").Value = '  -2.58%  '

$ws.Range("D20").Value = '314.09'
$ws.Range("E20").Value = '  -1.34%  '

$ws.Range("D21").Value = '4.05'
$ws.Range("E21").Value = '  -1.55%  '

$ws.Range("D22").Value = '6.47'
$ws.Range("E22").Value = '  -5.64%  '

$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").Value = '62.51'
$ws.Range("E24").Value = '  -1.96%  '

$ws.Range("D25").Value = '1.66'
$ws.Range("E25").Value = '  -9.94%  '

$ws.Range("D26").Value = '8.23'
$ws.Range("E26").Value = '  +5.86%  '

$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.43%  '

$ws.Range("D28").Value = '2.448.34'
$ws.Range("E28").Value = '  -4.63%  '

$ws.Range("D29").Value = '7.83'
$ws.Range("E29").Value = '  -4.30%  '

$ws.Range("D30").Value = '1.36'
$ws.Range("E30").Value = '  -8.11%  '

$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '489.00'
$ws.Range("E31").Value = '  -6.86%  '

$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = '0.0₃0851'
$ws.Range("E32").Value = '  -11.97%  '

$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").Value = '0.143'
$ws.Range("E33").Value = '  -2.26%  '

$ws.Range("D34").Value = '1.78'
$ws.Range("E34").Value = '  -5.05%  '

$ws.Range("E35").Value = '  -6.47%  '

$ws.Range("D36").Value = '0.992'
$ws.Range("E36").Value = '  -0.69%  '

$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '4.51'
$ws.Range("E37").Value = '  -4.19%  '

$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").Value = '0.371'
$ws.Range("E38").Value = '  -2.48%  '

$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").Value = '18.38'
$ws.Range("E39").Value = '  +1.09%  '

$ws.Range("D40").Value = '5.13'
$ws.Range("E40").Value = '  -8.85%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '1.75'
$ws.Range("E41").Value = '  +0.83%  '

$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = '141.19'
$ws.Range("E42").Value = '  +1.89%  '

$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").Value = '140.35'
$ws.Range("E45").Value = '  -1.64%  '

$ws.Range("D46").Value = '3.52'
$ws.Range("E46").Value = '  -2.30%  '

$ws.Range("D47").Value = '2.03'
$ws.Range("E47").Value = '  -10.33%  '

$ws.Range("D48").Value = '0.0506'
$ws.Range("E48").Value = '  -3.38%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '18.94'
$ws.Range("E49").Value = '  -10.53%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '0.563'
$ws.Range("E50").Value = '  -3.86%  '

$ws.Range("D51").Value = '0.0890'
$ws.Range("E51").Value = '  -4.16%  '
